$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "removeStuff" block (rows 38-39) mirrors the existing
# "removeDepartment" block (rows 26-27): a header row followed by a
# data row, same column layout (A name, B "1", C Company method,
# D/E key-value pair for rowIndex/isRemove).

# Row 38 (header row): copy formatting then values from row 26 (A:E)
$ws.Range("A26:E26").Copy() | Out-Null
$ws.Range("A38:E38").PasteSpecial(-4122) | Out-Null
$ws.Range("A26:E26").Copy() | Out-Null
$ws.Range("A38:E38").PasteSpecial(-4104) | Out-Null

# Row 39 (data row): copy formatting then values from row 27 (A:F,
# F27 is the trailing bordered blank cell)
$ws.Range("A27:F27").Copy() | Out-Null
$ws.Range("A39:F39").PasteSpecial(-4122) | Out-Null
$ws.Range("A27:F27").Copy() | Out-Null
$ws.Range("A39:F39").PasteSpecial(-4104) | Out-Null

# Rename the automation script + controller method for the new
# "remove_stuff" script.
$ws.Range("A38").Value = "TrainScheduling_ltrailways_removeStuff"
$ws.Range("C38").Value = "CompanyManagement.removeStuff"

$ws.Range("A39").Value = "TrainScheduling_ltrailways_removeStuff"
$ws.Range("C39").Value = "CompanyManagement.removeStuff"

$excel.CutCopyMode = 0

$ws.Activate() | Out-Null
$ws.Range("E39").Select() | Out-Null
